$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    if ($cell.Text -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
